$d = $word.ActiveDocument

# The document currently opens with:
#   Paragraph 1 (pStyle Heading1): "Three Who Have Died"
#   Paragraph 2 (no pStyle, bold run): "By Dorothy Day"
# wrapped by a bookmarkStart/bookmarkEnd pair around paragraph 1.
#
# Target (pandoc-style title block):
#   Paragraph 1 (pStyle Title): "Three" " " "Who" " " "Have" " " "Died"  (separate runs)
#   Paragraph 2 (pStyle Authors): "Dorothy" " " "Day"  (separate runs)
#
# Each paragraph is rewritten via InsertXML scoped to that single paragraph's
# own range (rather than one XML blob spanning both paragraphs). Doing the
# replacement one paragraph at a time keeps any pre-existing bookmark that
# anchors to paragraph 1 tightly wrapped around paragraph 1 alone, instead of
# letting it balloon out to also enclose paragraph 2.

$titleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:pStyle w:val="Title"/></w:pPr>
            <w:r><w:t xml:space="preserve">Three</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Who</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Have</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Died</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$null = $d.Paragraphs(1).Range.InsertXML($titleXml)

$authorsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:pStyle w:val="Authors"/></w:pPr>
            <w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$null = $d.Paragraphs(2).Range.InsertXML($authorsXml)
